$wb = $excel.ActiveWorkbook

# Rename sheet "ForgotPassword4A" to "ForgotPassword4B"
$ws = $wb.Worksheets.Item("ForgotPassword4A")
$ws.Name = "ForgotPassword4B"

# Update selection on that sheet to C15
$ws.Activate()
$ws.Range("C15").Select()
